$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 63, shifting rows 63-121 down to 64-122.
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new weekly record.
$ws.Cells.Item(63, 1).Value2 = 4
$ws.Cells.Item(63, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(63, 3).Value2 = "Los Lagos"
$ws.Cells.Item(63, 4).Value2 = 44483
$ws.Cells.Item(63, 5).Value2 = 10
$ws.Cells.Item(63, 6).Value2 = 100112028
$ws.Cells.Item(63, 7).Value2 = "Sandia"
$ws.Cells.Item(63, 8).Value2 = "Sin especificar"
$ws.Cells.Item(63, 9).Value2 = "Primera"
$ws.Cells.Item(63, 10).Value2 = 450
$ws.Cells.Item(63, 11).Value2 = 1000
$ws.Cells.Item(63, 12).Value2 = 1000
$ws.Cells.Item(63, 13).Value2 = 1000
$ws.Cells.Item(63, 14).Value2 = "$/kilo (volumen en unidades)"
$ws.Cells.Item(63, 15).Value2 = "Perú"
$ws.Cells.Item(63, 16).Value2 = 1000
$ws.Cells.Item(63, 17).Value2 = 1
$ws.Cells.Item(63, 18).Value2 = "Hortaliza"
